$wb = $excel.ActiveWorkbook

# --- Localization-status report regenerated for handoff ---
# Status flips from "In Translation" to "Ready for handoff" and the
# "generated at" timestamps advance a bit on every sheet that tracks them.
# The shared status text / timestamp text is literal (text-formatted dates,
# not real date serials), so set the cell .Value directly.

$newStatus = "Ready for handoff"

# Column width: the engine quantizes ColumnWidth (character units) onto a
# 1/6-character pixel grid before writing the OOXML <col width>. Picking
# 16.333333333333332 (=98/6) lands as close as this grid allows to the
# author's recorded width of 17.2159881591797.
$targetColWidth = 16.333333333333332

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-27 12:36:37"
$wsOverview.Columns.Item(5).ColumnWidth = $targetColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColWidth

# ----- zh-cn sheet -----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-27 12:36:33"
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColWidth

# ----- de-de sheet -----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-27 12:36:37"
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColWidth
